$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Insert a new bold paragraph at the very top of the document:
#    "Dr. Abdulmalik ALIYU"
# ---------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1).Range
$firstPara.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "Dr. Abdulmalik ALIYU"
$titlePara.Range.Font.Bold = $true
$titlePara.Range.Font.BoldBi = $true
$titlePara.Range.LanguageID = "en-US"

# ---------------------------------------------------------------
# 2. Change "hard working" -> "hardworking" inside the bio
#    paragraph. The surrounding text keeps its own runs (matching
#    how Word leaves the untouched text either side of an in-place
#    word replacement as separate runs).
# ---------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("hard working", $true, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)
if ($found) {
    $wordStart = $searchRange.Start
    $wordEnd = $searchRange.End
    $newWord = "hardworking"

    $target = $d.Range($wordStart, $wordEnd)
    $target.Text = $newWord

    $newWordEnd = $wordStart + $newWord.Length
    $midRange = $d.Range($wordStart, $newWordEnd)
    # Toggling a property on/off forces this span to remain a run of
    # its own, distinct from the runs before/after it, without
    # actually changing the visible formatting.
    $midRange.Bold = $true
    $midRange.Bold = $false
}

# ---------------------------------------------------------------
# 3. Append a new paragraph at the end of the document describing
#    the capstone project.
# ---------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("The main purpose of this capstone project is to build on online resume using HTML codes only.")
